$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values for rows 2-38 replacing the old computed values.
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 3
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 2
    14 = 6
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 0
    33 = 2
    34 = 2
    35 = 2
    36 = 1
    37 = 2
    38 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
